# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. This updates the DAMSLTag (column I) and
# DialogAct (column J) values for a set of rows in Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{Row=13; Tag="sd"; Act="Statement-non-opinion"},
    @{Row=40; Tag="sv"; Act="Statement-opinion"},
    @{Row=61; Tag="b"; Act="Acknowledge (Backchannel)"},
    @{Row=90; Tag="sv"; Act="Statement-opinion"},
    @{Row=91; Tag="sd"; Act="Statement-non-opinion"},
    @{Row=101; Tag="sv"; Act="Statement-opinion"},
    @{Row=120; Tag="%"; Act="Uninterpretable"},
    @{Row=127; Tag="sv"; Act="Statement-opinion"},
    @{Row=132; Tag="%"; Act="Uninterpretable"},
    @{Row=146; Tag="aa"; Act="Agree/Accept"},
    @{Row=147; Tag="aa"; Act="Agree/Accept"},
    @{Row=161; Tag="sv"; Act="Statement-opinion"},
    @{Row=192; Tag="aa"; Act="Agree/Accept"},
    @{Row=207; Tag="aa"; Act="Agree/Accept"},
    @{Row=213; Tag="sd"; Act="Statement-non-opinion"},
    @{Row=216; Tag="sd"; Act="Statement-non-opinion"},
    @{Row=229; Tag="sv"; Act="Statement-opinion"},
    @{Row=233; Tag="sd"; Act="Statement-non-opinion"},
    @{Row=239; Tag="sv"; Act="Statement-opinion"},
    @{Row=246; Tag="sv"; Act="Statement-opinion"},
    @{Row=258; Tag="%"; Act="Uninterpretable"},
    @{Row=266; Tag="b"; Act="Acknowledge (Backchannel)"},
    @{Row=272; Tag="sd"; Act="Statement-non-opinion"},
    @{Row=274; Tag="aa"; Act="Agree/Accept"},
    @{Row=279; Tag="aa"; Act="Agree/Accept"},
    @{Row=280; Tag="aa"; Act="Agree/Accept"},
    @{Row=295; Tag="b"; Act="Acknowledge (Backchannel)"},
    @{Row=300; Tag="sv"; Act="Statement-opinion"},
    @{Row=315; Tag="ba"; Act="Appreciation"},
    @{Row=320; Tag="%"; Act="Uninterpretable"},
    @{Row=329; Tag="sd"; Act="Statement-non-opinion"},
    @{Row=334; Tag="sd"; Act="Statement-non-opinion"},
    @{Row=336; Tag="sd"; Act="Statement-non-opinion"},
    @{Row=340; Tag="b"; Act="Acknowledge (Backchannel)"},
    @{Row=347; Tag="%"; Act="Uninterpretable"},
    @{Row=353; Tag="aa"; Act="Agree/Accept"},
    @{Row=354; Tag="%"; Act="Uninterpretable"},
    @{Row=373; Tag="sd"; Act="Statement-non-opinion"}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.Tag
    $ws.Cells.Item($u.Row, 10).Value = $u.Act
}
